$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 12 ("renault espace 2022"), a duplicate of row 8's data.
# This shifts rows 13-15 up to become rows 12-14.
$ws.Rows.Item(12).Delete()

# Set column A width (raw OOXML width=15; COM ColumnWidth adds 5/6 padding)
$ws.Columns.Item(1).ColumnWidth = 14.1666666666667

# Simplify car names in column A to brand names only
$ws.Range("A2").Value = "Ford"
$ws.Range("A3").Value = "ZEEKR"
$ws.Range("A4").Value = "MAXUS"
$ws.Range("A5").Value = "VW"
$ws.Range("A6").Value = "Škoda"
$ws.Range("A7").Value = "BMW"
$ws.Range("A8").Value = "Renault"
$ws.Range("A9").Value = "Mercedes-Benz"
$ws.Range("A10").Value = "Suzuki"
$ws.Range("A11").Value = "Dacia"
$ws.Range("A12").Value = "Toyota"
$ws.Range("A13").Value = "Honda"
$ws.Range("A14").Value = "NIO"
